$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1721.9474
$ws.Range("I15").Value = 1721.9474
$ws.Range("K15").Value = 5165.8422
$ws.Range("M15").Value = -4996.8422
$ws.Range("H51").Value = 4537.273
$ws.Range("I51").Value = 7500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -7016
$ws.Range("H53").Value = 1888.3529
$ws.Range("I53").Value = 72.5
$ws.Range("J53").Value = 2447.077
$ws.Range("K53").Value = 72.5
$ws.Range("L53").Value = 2447.077
$ws.Range("M53").Value = 564.5
$ws.Range("N53").Value = -3721.077
$ws.Range("H98").Value = 891.94116
$ws.Range("I98").Value = 914.4666999999999
$ws.Range("K98").Value = 914.4666999999999
$ws.Range("M98").Value = 583.5333000000001
$ws.Range("H100").Value = 66668372
$ws.Range("I100").Value = 142857870
$ws.Range("K100").Value = 142857870
$ws.Range("M100").Value = -142857329
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H106").Value = 11907045
$ws.Range("H122").Value = 891.94116
$ws.Range("I122").Value = 914.4666999999999
$ws.Range("K122").Value = 2743.4001
$ws.Range("M122").Value = -293.4000999999998
$ws.Range("H129").Value = 168127.7
$ws.Range("J129").Value = 176951.98
$ws.Range("L129").Value = 530855.9400000001
$ws.Range("N129").Value = -540855.9400000001
$ws.Range("H132").Value = 6740.5835
$ws.Range("I132").Value = 8654.111000000001
$ws.Range("K132").Value = 25962.333
$ws.Range("M132").Value = -23432.333
$ws.Range("H137").Value = 25720.072
$ws.Range("I137").Value = 1266.7097
$ws.Range("K137").Value = 3800.1291
$ws.Range("M137").Value = -1250.1291
$ws.Range("H141").Value = 1520.2941
$ws.Range("I141").Value = 1365.3125
$ws.Range("K141").Value = 4095.9375
$ws.Range("M141").Value = 1084.0625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22879.893
$ws.Range("I32").Value = 25948.584
$ws.Range("K32").Value = 25948.584
$ws.Range("M32").Value = -25661.584
$ws.Range("H88").Value = 1000054
$ws.Range("J88").Value = 1000054
$ws.Range("L88").Value = 1000054
$ws.Range("N88").Value = -1000866
$ws.Range("H91").Value = 1000054
$ws.Range("J91").Value = 1000054
$ws.Range("L91").Value = 1000054
$ws.Range("N91").Value = -1002862

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1555.421
$ws.Range("J99").Value = 1602.75
$ws.Range("L99").Value = 1602.75
$ws.Range("N99").Value = -4598.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4729.7144
$ws.Range("J4").Value = 4729.7144
$ws.Range("L4").Value = 4729.7144
$ws.Range("N4").Value = -4953.7144
$ws.Range("H31").Value = 14657.223
$ws.Range("I31").Value = 54996.668
$ws.Range("J31").Value = 3131.6667
$ws.Range("K31").Value = 54996.668
$ws.Range("L31").Value = 3131.6667
$ws.Range("M31").Value = -54701.668
$ws.Range("N31").Value = -3721.6667
$ws.Range("H34").Value = 14657.223
$ws.Range("I34").Value = 54996.668
$ws.Range("J34").Value = 3131.6667
$ws.Range("K34").Value = 54996.668
$ws.Range("L34").Value = 3131.6667
$ws.Range("M34").Value = -54794.668
$ws.Range("N34").Value = -3535.6667
$ws.Range("H68").Value = 34067.6
$ws.Range("J68").Value = 34067.6
$ws.Range("L68").Value = 34067.6
$ws.Range("N68").Value = -35565.6
$ws.Range("H71").Value = 34067.6
$ws.Range("J71").Value = 34067.6
$ws.Range("L71").Value = 102202.8
$ws.Range("N71").Value = -109690.8
$ws.Range("H105").Value = 25000610
$ws.Range("I105").Value = 25000610
$ws.Range("K105").Value = 25000610
$ws.Range("M105").Value = -24998863
$ws.Range("H107").Value = 1135.75
$ws.Range("I107").Value = 380.75
$ws.Range("J107").Value = 1513.25
$ws.Range("K107").Value = 380.75
$ws.Range("L107").Value = 1513.25
$ws.Range("M107").Value = 1539.25
$ws.Range("N107").Value = -5353.25
$ws.Range("H132").Value = 21526.076
$ws.Range("I132").Value = 22820
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 68460
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -65930
$ws.Range("N132").Value = -23057

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 248
$ws.Range("I8").Value = 248
$ws.Range("K8").Value = 744
$ws.Range("M8").Value = -605
$ws.Range("H10").Value = 425.6
$ws.Range("I10").Value = 383.77777
$ws.Range("J10").Value = 802
$ws.Range("K10").Value = 1151.33331
$ws.Range("L10").Value = 2406
$ws.Range("M10").Value = -1012.33331
$ws.Range("N10").Value = -2684
$ws.Range("H61").Value = 180
$ws.Range("I61").Value = 50
$ws.Range("J61").Value = 266.66666
$ws.Range("K61").Value = 150
$ws.Range("L61").Value = 799.9999799999999
$ws.Range("M61").Value = 65
$ws.Range("N61").Value = -1229.99998
$ws.Range("H118").Value = 41668760
$ws.Range("I118").Value = 83333540
$ws.Range("K118").Value = 250000620
$ws.Range("M118").Value = -249999377
$ws.Range("H129").Value = 228159
$ws.Range("I129").Value = 861.1111
$ws.Range("J129").Value = 385519.06
$ws.Range("K129").Value = 2583.3333
$ws.Range("L129").Value = 1156557.18
$ws.Range("M129").Value = 2416.6667
$ws.Range("N129").Value = -1166557.18
$ws.Range("H131").Value = 817.73
$ws.Range("J131").Value = 817.73
$ws.Range("L131").Value = 2453.19
$ws.Range("N131").Value = -12533.19

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6500
$ws.Range("I5").Value = 3000
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -2888
$ws.Range("N5").Value = -10224
$ws.Range("H97").Value = 629.3333
$ws.Range("I97").Value = 629.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 629.3333
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -133.3333
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 35715750
$ws.Range("I102").Value = 38462930
$ws.Range("K102").Value = 38462930
$ws.Range("M102").Value = -38461308

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4919.5713
$ws.Range("I7").Value = 2841.3
$ws.Range("K7").Value = 2841.3
$ws.Range("M7").Value = -2729.3
$ws.Range("H40").Value = 4944.5
$ws.Range("I40").Value = 3751.5
$ws.Range("J40").Value = 6137.5
$ws.Range("K40").Value = 3751.5
$ws.Range("L40").Value = 6137.5
$ws.Range("M40").Value = -3615.5
$ws.Range("N40").Value = -6409.5
$ws.Range("H126").Value = 4919.5713
$ws.Range("I126").Value = 2841.3
$ws.Range("K126").Value = 8523.900000000001
$ws.Range("M126").Value = -6053.900000000001
$ws.Range("H136").Value = 16427.281
$ws.Range("I136").Value = 19951.076
$ws.Range("J136").Value = 1157.5
$ws.Range("K136").Value = 59853.228
$ws.Range("L136").Value = 3472.5
$ws.Range("M136").Value = -57303.228
$ws.Range("N136").Value = -8572.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1316.1538
$ws.Range("I122").Value = 1182.8182
$ws.Range("J122").Value = 2049.5
$ws.Range("K122").Value = 3548.4546
$ws.Range("L122").Value = 6148.5
$ws.Range("M122").Value = -1098.4546
$ws.Range("N122").Value = -11048.5
$ws.Range("H126").Value = 1895.7858
$ws.Range("I126").Value = 1129.2
$ws.Range("K126").Value = 3387.6
$ws.Range("M126").Value = -917.6000000000004
$ws.Range("H132").Value = 3140.4
$ws.Range("I132").Value = 1901.3334
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5704.0002
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3174.0002
$ws.Range("N132").Value = -20057
